# fichamentoTcc2019 — remove the two leftover "pg 50." / "pg 80" reference
# notes from the citation blocks (rows 7 and 15) but keep their existing
# cell formatting (borders / centered style) untouched, then move the
# window's active selection down to the now-empty E15 cell, scrolled so
# row 7 is visible at the top — matching where the author was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear just the text content of E7 ("pg 50.") and E15 ("pg 80"); the
# cells keep their original style (border + center alignment) because
# ClearContents only touches the value, not the formatting.
$ws.Range("E7").ClearContents()
$ws.Range("E15").ClearContents()

# Scroll the window so row 7 is the first visible row and move the
# selection to the single cell E15 (previously the selection was
# G7:H7 with G7 active).
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 7
    $win.ScrollColumn = 1
}
$ws.Range("E15").Select()
